# New API Query - 2023 Included
# API query to UN performed 11/26/2023. Query modified to include 2023 data.
#
# Updates on sheet "fromCSV":
#   - B2 (short-url): regenerate the API query short-url "ooGnM5" -> "q3XUxt"
#   - U2 (oip): blank/no-data marker "null" -> "-"
#   - V2 (hst): blank/no-data marker 0 -> "-", now left-aligned like the
#     other text placeholders (matches U2's formatting) instead of the
#     right-aligned numeric style it had before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "q3XUxt"
$ws.Range("U2").Value = "-"
$ws.Range("V2").Value = "-"
$ws.Range("V2").HorizontalAlignment = -4131
